$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'300.40"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'2.05%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'32.37"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'3.93%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'1.28%"
$ws.Range("E4").Style = "Normal"
$ws.Range("E5").Value = "'5.20%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'2.285"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'0.11%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'7.932"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'2.41%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.9226"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'1.53%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.09982"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'23.91%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1760"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'4.06%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.08427"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'4.52%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.03305"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'6.61%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.09863"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-2.13%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.001480"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-2.53%"
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'-0.74%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'3.515"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'1.09%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.821"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'1.93%"
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'5.14%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.3356"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'0.84%"
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'2.49%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'4.377"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'10.17%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.2087"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-0.71%"
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'0.01%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001215"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'0.37%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004368"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-2.81%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D27").Value = "'0.0003373"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'-0.75%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = "'0.01706"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'6.58%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.04680"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'5.30%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007724"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'4.88%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.009769"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'13.21%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.1392"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'4.70%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.002080"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'6.75%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.009612"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'0.91%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006068"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'1.48%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000745"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-0.72%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'2.794"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'24.68%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.001987"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-31.44%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002087"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'-0.72%"
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'-0.72%"
$ws.Range("E51").Style = "Normal"
